$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look numeric stay as text (matching source formatting),
# mirroring the original inline-string cell type for the Price column.
$textCells = @("D5","D6","D7","D10","D11","D12","D14","D16","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D41","D43","D44","D45","D46","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value2 = "41.185.43"
$ws.Range("E2").Value2 = "  -1.45%  "

# Row 3
$ws.Range("D3").Value2 = "2.178.50"
$ws.Range("E3").Value2 = "  -2.13%  "

# Row 4
$ws.Range("E4").Value2 = "  -0.11%  "

# Row 5
$ws.Range("D5").Value2 = "251.59"
$ws.Range("E5").Value2 = "  +0.26%  "

# Row 6
$ws.Range("D6").Value2 = "0.614"
$ws.Range("E6").Value2 = "  -2.46%  "

# Row 7
$ws.Range("D7").Value2 = "66.43"
$ws.Range("E7").Value2 = "  -7.73%  "

# Row 8
$ws.Range("E8").Value2 = "  -0.01%  "

# Row 9
$ws.Range("E9").Value2 = "  -2.50%  "

# Row 10
$ws.Range("D10").Value2 = "59.06"
$ws.Range("E10").Value2 = "  +1.45%  "

# Row 11
$ws.Range("D11").Value2 = "36.30"
$ws.Range("E11").Value2 = "  -12.05%  "

# Row 12
$ws.Range("D12").Value2 = "0.0935"
$ws.Range("E12").Value2 = "  -3.59%  "

# Row 13
$ws.Range("E13").Value2 = "  -1.47%  "

# Row 14
$ws.Range("D14").Value2 = "6.85"
$ws.Range("E14").Value2 = "  -4.38%  "

# Row 15
$ws.Range("D15").Value2 = "2.504.83"
$ws.Range("E15").Value2 = "  -2.06%  "

# Row 16
$ws.Range("D16").Value2 = "14.34"
$ws.Range("E16").Value2 = "  -4.42%  "

# Row 17
$ws.Range("E17").Value2 = "  -2.03%  "

# Row 18
$ws.Range("D18").Value2 = "2.181.01"
$ws.Range("E18").Value2 = "  -2.21%  "

# Row 19
$ws.Range("D19").Value2 = "41.119.24"
$ws.Range("E19").Value2 = "  -1.52%  "

# Row 20
$ws.Range("E20").Value2 = "  -1.77%  "

# Row 21
$ws.Range("E21").Value2 = "  -1.65%  "

# Row 22
$ws.Range("E22").Value2 = "  -2.41%  "

# Row 23
$ws.Range("D23").Value2 = "230.56"
$ws.Range("E23").Value2 = "  -2.07%  "

# Row 24
$ws.Range("D24").Value2 = "2.04"
$ws.Range("E24").Value2 = "  -3.61%  "

# Row 25
$ws.Range("D25").Value2 = "3.80"
$ws.Range("E25").Value2 = "  -11.14%  "

# Row 26
$ws.Range("D26").Value2 = "11.49"
$ws.Range("E26").Value2 = "  +7.14%  "

# Row 27
$ws.Range("E27").Value2 = "  +0.10%  "

# Row 28
$ws.Range("D28").Value2 = "2.42"
$ws.Range("E28").Value2 = "  -4.66%  "

# Row 29
$ws.Range("B29").Value2 = "Monero"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value2 = "168.46"
$ws.Range("E29").Value2 = "  -1.69%  "

# Row 30
$ws.Range("B30").Value2 = "Toncoin"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value2 = "2.03"
$ws.Range("E30").Value2 = "  -7.29%  "

# Row 31
$ws.Range("D31").Value2 = "20.24"
$ws.Range("E31").Value2 = "  -2.57%  "

# Row 32
$ws.Range("D32").Value2 = "0.123"
$ws.Range("E32").Value2 = "  -1.85%  "

# Row 33
$ws.Range("D33").Value2 = "5.76"
$ws.Range("E33").Value2 = "  +3.03%  "

# Row 34
$ws.Range("D34").Value2 = "0.0753"
$ws.Range("E34").Value2 = "  +2.66%  "

# Row 35
$ws.Range("E35").Value2 = "  -3.31%  "

# Row 36
$ws.Range("D36").Value2 = "4.53"
$ws.Range("E36").Value2 = "  -4.12%  "

# Row 37
$ws.Range("D37").Value2 = "3.94"
$ws.Range("E37").Value2 = "  -1.55%  "

# Row 38
$ws.Range("D38").Value2 = "24.53"
$ws.Range("E38").Value2 = "  -7.07%  "

# Row 39
$ws.Range("E39").Value2 = "  -0.11%  "

# Row 40
$ws.Range("E40").Value2 = "  -3.14%  "

# Row 41
$ws.Range("D41").Value2 = "5.33"
$ws.Range("E41").Value2 = "  +8.24%  "

# Row 42
$ws.Range("E42").Value2 = "  -7.72%  "

# Row 43
$ws.Range("D43").Value2 = "11.39"
$ws.Range("E43").Value2 = "  -5.64%  "

# Row 44
$ws.Range("D44").Value2 = "60.71"
$ws.Range("E44").Value2 = "  -9.25%  "

# Row 45
$ws.Range("D45").Value2 = "8.52"
$ws.Range("E45").Value2 = "  -3.07%  "

# Row 46
$ws.Range("D46").Value2 = "0.0997"
$ws.Range("E46").Value2 = "  -2.43%  "

# Row 47
$ws.Range("D47").Value2 = "0.189"
$ws.Range("E47").Value2 = "  -7.05%  "

# Row 48
$ws.Range("E48").Value2 = "  -0.07%  "

# Row 49
$ws.Range("E49").Value2 = "  -1.78%  "

# Row 50
$ws.Range("E50").Value2 = "  -4.31%  "

# Row 51
$ws.Range("B51").Value2 = "SynthetixNetwork"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value2 = "4.17"
$ws.Range("E51").Value2 = "  -10.15%  "
